$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("qemu_ipc")
$ws.Range("M1:AA29").Cut($ws.Range("AM14"))
